{"js": "// Append two new paragraphs to the end of the document body, after the\n// existing \"Prueba 6. Nuevos cambios 8:58 03-03-2020\" paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst p7 = lastParagraph.insertParagraph(\n  \"Prueba 7. Nuevos cambios: 9:05 03-03-2020\",\n  Word.InsertLocation.after\n);\nconst p8 = p7.insertParagraph(\n  \"Prueba 8: Nuevos cambios: 9:07 03-03-2020.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Append two new paragraphs to the end of the document, after the\n# existing \"Prueba 6. Nuevos cambios 8:58 03-03-2020\" paragraph.\n$d = $word.ActiveDocument\n\n$lastPara = $d.Paragraphs.Last\n$r = $lastPara.Range\n$r.Collapse(0)  # wdCollapseEnd\n\n$r.InsertParagraphAfter()\n$r = $d.Paragraphs.Last.Range\n$r.InsertAfter(\"Prueba 7. Nuevos cambios: 9:05 03-03-2020\")\n\n$r = $d.Paragraphs.Last.Range\n$r.InsertParagraphAfter()\n$r = $d.Paragraphs.Last.Range\n$r.InsertAfter(\"Prueba 8: Nuevos cambios: 9:07 03-03-2020.\")\n"}
